$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 5114.2383  # H64: 5174.95 -> 5114.2383
$ws.Cells.Item(64, 9).Value = 5421.4287  # I64: 5538.4614 -> 5421.4287
$ws.Cells.Item(64, 11).Value = 5421.4287  # K64: 5538.4614 -> 5421.4287
$ws.Cells.Item(64, 13).Value = -5173.4287  # M64: -5290.4614 -> -5173.4287
$ws.Cells.Item(67, 8).Value = 5114.2383  # H67: 5174.95 -> 5114.2383
$ws.Cells.Item(67, 9).Value = 5421.4287  # I67: 5538.4614 -> 5421.4287
$ws.Cells.Item(67, 11).Value = 5421.4287  # K67: 5538.4614 -> 5421.4287
$ws.Cells.Item(67, 13).Value = -4563.4287  # M67: -4680.4614 -> -4563.4287
$ws.Cells.Item(76, 8).Value = 5324.25  # H76: 5202.533 -> 5324.25
$ws.Cells.Item(76, 9).Value = 5224.4165  # I76: 5003.5835 -> 5224.4165
$ws.Cells.Item(76, 10).Value = 5623.75  # J76: 5998.3335 -> 5623.75
$ws.Cells.Item(76, 11).Value = 5224.4165  # K76: 5003.5835 -> 5224.4165
$ws.Cells.Item(76, 12).Value = 5623.75  # L76: 5998.3335 -> 5623.75
$ws.Cells.Item(76, 13).Value = -4909.4165  # M76: -4688.5835 -> -4909.4165
$ws.Cells.Item(76, 14).Value = -6253.75  # N76: -6628.3335 -> -6253.75
$ws.Cells.Item(79, 8).Value = 5324.25  # H79: 5202.533 -> 5324.25
$ws.Cells.Item(79, 9).Value = 5224.4165  # I79: 5003.5835 -> 5224.4165
$ws.Cells.Item(79, 10).Value = 5623.75  # J79: 5998.3335 -> 5623.75
$ws.Cells.Item(79, 11).Value = 5224.4165  # K79: 5003.5835 -> 5224.4165
$ws.Cells.Item(79, 12).Value = 5623.75  # L79: 5998.3335 -> 5623.75
$ws.Cells.Item(79, 13).Value = -4132.4165  # M79: -3911.5835 -> -4132.4165
$ws.Cells.Item(79, 14).Value = -7807.75  # N79: -8182.3335 -> -7807.75
$ws.Cells.Item(129, 8).Value = 2350.4656  # H129: 2263.3386 -> 2350.4656
$ws.Cells.Item(129, 9).Value = 1039.75  # I129: 1019.875 -> 1039.75
$ws.Cells.Item(129, 11).Value = 3119.25  # K129: 3059.625 -> 3119.25
$ws.Cells.Item(129, 13).Value = 1880.75  # M129: 1940.375 -> 1880.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2099.7896  # H2: 2084.6667 -> 2099.7896
$ws.Cells.Item(2, 9).Value = 1684.6666  # I2: 1677.68 -> 1684.6666
$ws.Cells.Item(2, 11).Value = 1684.6666  # K2: 1677.68 -> 1684.6666
$ws.Cells.Item(2, 13).Value = -1571.6666  # M2: -1564.68 -> -1571.6666
$ws.Cells.Item(4, 8).Value = 672.75  # H4: 1045 -> 672.75
$ws.Cells.Item(4, 9).Value = 363.66666  # I4: 490 -> 363.66666
$ws.Cells.Item(4, 11).Value = 363.66666  # K4: 490 -> 363.66666
$ws.Cells.Item(4, 13).Value = -247.66666  # M4: -374 -> -247.66666
$ws.Cells.Item(15, 8).Value = 9759.6  # H15: 13299.667 -> 9759.6
$ws.Cells.Item(15, 10).Value = 9759.6  # J15: 13299.667 -> 9759.6
$ws.Cells.Item(15, 12).Value = 9759.6  # L15: 13299.667 -> 9759.6
$ws.Cells.Item(15, 14).Value = -10459.6  # N15: -13999.667 -> -10459.6
$ws.Cells.Item(61, 8).Value = 2023.0555  # H61: 1912.7646 -> 2023.0555
$ws.Cells.Item(61, 10).Value = 3899  # J61: 3899.5 -> 3899
$ws.Cells.Item(61, 12).Value = 3899  # L61: 3899.5 -> 3899
$ws.Cells.Item(61, 14).Value = -4323  # N61: -4323.5 -> -4323
$ws.Cells.Item(116, 8).Value = 2099.7896  # H116: 2084.6667 -> 2099.7896
$ws.Cells.Item(116, 9).Value = 1684.6666  # I116: 1677.68 -> 1684.6666
$ws.Cells.Item(116, 11).Value = 1684.6666  # K116: 1677.68 -> 1684.6666
$ws.Cells.Item(116, 13).Value = 609.3334  # M116: 616.3199999999999 -> 609.3334
$ws.Cells.Item(132, 8).Value = 18521000  # H132: 19233338 -> 18521000
$ws.Cells.Item(132, 9).Value = 2519.6  # I132: 2641.6843 -> 2519.6
$ws.Cells.Item(132, 11).Value = 7558.799999999999  # K132: 7925.0529 -> 7558.799999999999
$ws.Cells.Item(132, 13).Value = -5028.799999999999  # M132: -5395.0529 -> -5028.799999999999
$ws.Cells.Item(136, 8).Value = 2023.0555  # H136: 1912.7646 -> 2023.0555
$ws.Cells.Item(136, 10).Value = 3899  # J136: 3899.5 -> 3899
$ws.Cells.Item(136, 12).Value = 11697  # L136: 11698.5 -> 11697
$ws.Cells.Item(136, 14).Value = -16797  # N136: -16798.5 -> -16797

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2099.7896  # H3: 2084.6667 -> 2099.7896
$ws.Cells.Item(3, 9).Value = 1684.6666  # I3: 1677.68 -> 1684.6666
$ws.Cells.Item(3, 11).Value = 1684.6666  # K3: 1677.68 -> 1684.6666
$ws.Cells.Item(3, 13).Value = -1570.6666  # M3: -1563.68 -> -1570.6666
$ws.Cells.Item(94, 8).Value = 13214.85  # H94: 13867.211 -> 13214.85
$ws.Cells.Item(94, 9).Value = 3961.4614  # I94: 4223.25 -> 3961.4614
$ws.Cells.Item(94, 11).Value = 3961.4614  # K94: 4223.25 -> 3961.4614
$ws.Cells.Item(94, 13).Value = -3510.4614  # M94: -3772.25 -> -3510.4614
$ws.Cells.Item(134, 8).Value = 15354320  # H134: 16207320 -> 15354320
$ws.Cells.Item(134, 9).Value = 7356576.5  # I134: 7816343 -> 7356576.5
$ws.Cells.Item(134, 11).Value = 22069729.5  # K134: 23449029 -> 22069729.5
$ws.Cells.Item(134, 13).Value = -22067194.5  # M134: -23446494 -> -22067194.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13, 8).Value = 1844.75  # H13: 1874.4546 -> 1844.75
$ws.Cells.Item(13, 9).Value = 2206  # I13: 2700 -> 2206
$ws.Cells.Item(13, 10).Value = 1586.7142  # J13: 1402.7142 -> 1586.7142
$ws.Cells.Item(13, 11).Value = 2206  # K13: 2700 -> 2206
$ws.Cells.Item(13, 12).Value = 1586.7142  # L13: 1402.7142 -> 1586.7142
$ws.Cells.Item(13, 13).Value = -2067  # M13: -2561 -> -2067
$ws.Cells.Item(13, 14).Value = -1864.7142  # N13: -1680.7142 -> -1864.7142
$ws.Cells.Item(14, 8).Value = 3800  # H14: 2399 -> 3800
$ws.Cells.Item(14, 9).Value = 0  # I14: 998 -> 0
$ws.Cells.Item(14, 11).Value = 0  # K14: 998 -> 0
$ws.Cells.Item(14, 13).ClearContents()  # M14: was -828
$ws.Cells.Item(25, 8).Value = 5000  # H25: 4670.3335 -> 5000
$ws.Cells.Item(25, 9).Value = 5000  # I25: 4670.3335 -> 5000
$ws.Cells.Item(25, 11).Value = 5000  # K25: 4670.3335 -> 5000
$ws.Cells.Item(25, 13).Value = -4826  # M25: -4496.3335 -> -4826
$ws.Cells.Item(26, 8).Value = 20000  # H26: 0 -> 20000
$ws.Cells.Item(26, 10).Value = 20000  # J26: 0 -> 20000
$ws.Cells.Item(26, 12).Value = 20000  # L26: 0 -> 20000
$ws.Cells.Item(26, 14).Value = -20574  # N26: None -> -20574
$ws.Cells.Item(31, 8).Value = 2044.9  # H31: 2138 -> 2044.9
$ws.Cells.Item(31, 10).Value = 2075.3845  # J31: 2212.8572 -> 2075.3845
$ws.Cells.Item(31, 12).Value = 2075.3845  # L31: 2212.8572 -> 2075.3845
$ws.Cells.Item(31, 14).Value = -2665.3845  # N31: -2802.8572 -> -2665.3845
$ws.Cells.Item(34, 8).Value = 2044.9  # H34: 2138 -> 2044.9
$ws.Cells.Item(34, 10).Value = 2075.3845  # J34: 2212.8572 -> 2075.3845
$ws.Cells.Item(34, 12).Value = 2075.3845  # L34: 2212.8572 -> 2075.3845
$ws.Cells.Item(34, 14).Value = -2479.3845  # N34: -2616.8572 -> -2479.3845
$ws.Cells.Item(37, 8).Value = 7000.2  # H37: 8000.25 -> 7000.2
$ws.Cells.Item(37, 10).Value = 6250  # J37: 7333.3335 -> 6250
$ws.Cells.Item(37, 12).Value = 6250  # L37: 7333.3335 -> 6250
$ws.Cells.Item(37, 14).Value = -6464  # N37: -7547.3335 -> -6464
$ws.Cells.Item(58, 8).Value = 2777.7273  # H58: 2881.9524 -> 2777.7273
$ws.Cells.Item(58, 9).Value = 1364.2727  # I58: 1409.7273 -> 1364.2727
$ws.Cells.Item(58, 10).Value = 4191.1816  # J58: 4501.4 -> 4191.1816
$ws.Cells.Item(58, 11).Value = 1364.2727  # K58: 1409.7273 -> 1364.2727
$ws.Cells.Item(58, 12).Value = 4191.1816  # L58: 4501.4 -> 4191.1816
$ws.Cells.Item(58, 13).Value = -1161.2727  # M58: -1206.7273 -> -1161.2727
$ws.Cells.Item(58, 14).Value = -4597.1816  # N58: -4907.4 -> -4597.1816
$ws.Cells.Item(132, 8).Value = 3065.3044  # H132: 2975.0833 -> 3065.3044
$ws.Cells.Item(132, 9).Value = 2798.6  # I132: 2798.7 -> 2798.6
$ws.Cells.Item(132, 10).Value = 4843.3335  # J132: 3857 -> 4843.3335
$ws.Cells.Item(132, 11).Value = 8395.799999999999  # K132: 8396.099999999999 -> 8395.799999999999
$ws.Cells.Item(132, 12).Value = 14530.0005  # L132: 11571 -> 14530.0005
$ws.Cells.Item(132, 13).Value = -5865.799999999999  # M132: -5866.099999999999 -> -5865.799999999999
$ws.Cells.Item(132, 14).Value = -19590.0005  # N132: -16631 -> -19590.0005
$ws.Cells.Item(134, 8).Value = 6251685.5  # H134: 6251884.5 -> 6251685.5
$ws.Cells.Item(134, 9).Value = 1797.9333  # I134: 2010 -> 1797.9333
$ws.Cells.Item(134, 11).Value = 5393.7999  # K134: 6030 -> 5393.7999
$ws.Cells.Item(134, 13).Value = -2858.7999  # M134: -3495 -> -2858.7999
$ws.Cells.Item(136, 8).Value = 2777.7273  # H136: 2881.9524 -> 2777.7273
$ws.Cells.Item(136, 9).Value = 1364.2727  # I136: 1409.7273 -> 1364.2727
$ws.Cells.Item(136, 10).Value = 4191.1816  # J136: 4501.4 -> 4191.1816
$ws.Cells.Item(136, 11).Value = 4092.8181  # K136: 4229.1819 -> 4092.8181
$ws.Cells.Item(136, 12).Value = 12573.5448  # L136: 13504.2 -> 12573.5448
$ws.Cells.Item(136, 13).Value = -1542.8181  # M136: -1679.1819 -> -1542.8181
$ws.Cells.Item(136, 14).Value = -17673.5448  # N136: -18604.2 -> -17673.5448

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 626  # H23: 410.83334 -> 626
$ws.Cells.Item(23, 9).Value = 500  # I23: 154.33333 -> 500
$ws.Cells.Item(23, 10).Value = 668  # J23: 667.3333 -> 668
$ws.Cells.Item(23, 11).Value = 1500  # K23: 462.99999 -> 1500
$ws.Cells.Item(23, 12).Value = 2004  # L23: 2001.9999 -> 2004
$ws.Cells.Item(23, 13).Value = -1265  # M23: -227.99999 -> -1265
$ws.Cells.Item(23, 14).Value = -2474  # N23: -2471.9999 -> -2474
$ws.Cells.Item(46, 8).Value = 2004390  # H46: 2505250 -> 2004390
$ws.Cells.Item(46, 9).Value = 5000475  # I46: 10000000 -> 5000475
$ws.Cells.Item(46, 11).Value = 15001425  # K46: 30000000 -> 15001425
$ws.Cells.Item(46, 13).Value = -15001334  # M46: -29999909 -> -15001334

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 2333  # H6: 2499.5 -> 2333
$ws.Cells.Item(6, 9).Value = 2500  # I6: 3000 -> 2500
$ws.Cells.Item(6, 11).Value = 2500  # K6: 3000 -> 2500
$ws.Cells.Item(6, 13).Value = -2387  # M6: -2887 -> -2387
$ws.Cells.Item(16, 8).Value = 2333  # H16: 2499.5 -> 2333
$ws.Cells.Item(16, 9).Value = 2500  # I16: 3000 -> 2500
$ws.Cells.Item(16, 11).Value = 2500  # K16: 3000 -> 2500
$ws.Cells.Item(16, 13).Value = -2250  # M16: -2750 -> -2250
$ws.Cells.Item(17, 8).Value = 2150.75  # H17: 1525.5 -> 2150.75
$ws.Cells.Item(17, 9).Value = 0  # I17: 250 -> 0
$ws.Cells.Item(17, 10).Value = 2150.75  # J17: 1780.6 -> 2150.75
$ws.Cells.Item(17, 11).Value = 0  # K17: 250 -> 0
$ws.Cells.Item(17, 12).Value = 2150.75  # L17: 1780.6 -> 2150.75
$ws.Cells.Item(17, 13).ClearContents()  # M17: was -82
$ws.Cells.Item(17, 14).Value = -2486.75  # N17: -2116.6 -> -2486.75
$ws.Cells.Item(113, 8).Value = 3302.0557  # H113: 3320.611 -> 3302.0557
$ws.Cells.Item(113, 9).Value = 3364.7144  # I113: 3388.5715 -> 3364.7144
$ws.Cells.Item(113, 11).Value = 3364.7144  # K113: 3388.5715 -> 3364.7144
$ws.Cells.Item(113, 13).Value = -1194.7144  # M113: -1218.5715 -> -1194.7144

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 6999.5  # H25: 7000 -> 6999.5
$ws.Cells.Item(25, 9).Value = 6999.5  # I25: 7000 -> 6999.5
$ws.Cells.Item(25, 11).Value = 6999.5  # K25: 7000 -> 6999.5
$ws.Cells.Item(25, 13).Value = -6769.5  # M25: -6770 -> -6769.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 0  # H24: 1995 -> 0
$ws.Cells.Item(24, 10).Value = 0  # J24: 1995 -> 0
$ws.Cells.Item(24, 12).Value = 0  # L24: 1995 -> 0
$ws.Cells.Item(24, 14).ClearContents()  # N24: was -2455
$ws.Cells.Item(132, 8).Value = 2996.0938  # H132: 2989.4062 -> 2996.0938
$ws.Cells.Item(132, 9).Value = 3074.1667  # I132: 3067.0334 -> 3074.1667
$ws.Cells.Item(132, 11).Value = 9222.500100000001  # K132: 9201.100199999999 -> 9222.500100000001
$ws.Cells.Item(132, 13).Value = -6692.500100000001  # M132: -6671.100199999999 -> -6692.500100000001
$ws.Cells.Item(136, 8).Value = 972.5  # H136: 962.64105 -> 972.5
$ws.Cells.Item(136, 9).Value = 941.4194  # I136: 956.1667 -> 941.4194
$ws.Cells.Item(136, 10).Value = 1110.1428  # J136: 984.2222 -> 1110.1428
$ws.Cells.Item(136, 11).Value = 2824.2582  # K136: 2868.5001 -> 2824.2582
$ws.Cells.Item(136, 12).Value = 3330.4284  # L136: 2952.6666 -> 3330.4284
$ws.Cells.Item(136, 13).Value = -274.2582000000002  # M136: -318.5001000000002 -> -274.2582000000002
$ws.Cells.Item(136, 14).Value = -8430.428400000001  # N136: -8052.6666 -> -8430.428400000001
